$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table (rows 3-6): update existing rows ---
# Row 3: was 22.110.1.1 -> becomes 22.240.0.6 entry
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.240.0.6"
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 1667
$ws.Range("D3").Value = 98.3

# Row 4: was 22.240.0.6 -> becomes 22.110.1.1 entry
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.110.1.1"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 351
$ws.Range("D4").Value = 98.40000000000001

# Row 5: Critical Minutes changes 689 -> 684
$ws.Range("C5").Value = 684

# Row 6: Critical Minutes changes 132 -> 134
$ws.Range("C6").Value = 134

# --- Insert a new row for a new "Bad Driver" entry (23.40.0.4) ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4"
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 573
$ws.Range("D7").Value = 98.90000000000001

# --- Totals row (now row 8 after the insert) ---
$ws.Range("B8").Value = 26
$ws.Range("C8").Value = 3409

Write-Host "done"
